# edit.ps1
# Applies the cryptos.xlsx data refresh described by the commit
# "Updated cryptos list on Fri Apr 19 13:45:41 UTC 2024 with GitHub Actions"
#
# The worksheet stores every data cell (columns B-E) as inline text strings.
# Several Price values look numeric (e.g. "1.00", "0.999") - if Excel is
# allowed to auto-detect these as numbers it will silently drop the text
# formatting (trailing zeros, thousand-dot separators, etc.), so every cell
# is forced to Text format before the value is written, then the cell style
# is reset back to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '64.657.49'
    'E2' = '  +4.06%  '
    'D3' = '3.094.52'
    'E3' = '  +2.50%  '
    'D4' = '1.00'
    'E4' = '  -0.13%  '
    'D5' = '559.60'
    'E5' = '  +3.40%  '
    'D6' = '143.68'
    'E6' = '  +8.46%  '
    'E7' = '  -0.10%  '
    'D8' = '3.093.27'
    'E8' = '  +2.65%  '
    'D9' = '0.500'
    'E9' = '  +1.61%  '
    'D10' = '6.32'
    'E10' = '  +3.37%  '
    'D11' = '0.152'
    'E11' = '  +4.22%  '
    'D12' = '0.473'
    'E12' = '  +6.72%  '
    'D13' = '0.0000229'
    'E13' = '  +3.74%  '
    'D14' = '35.29'
    'E14' = '  +4.07%  '
    'D15' = '3.598.86'
    'E15' = '  +2.49%  '
    'D16' = '64.758.09'
    'E16' = '  +4.09%  '
    'D17' = '3.099.71'
    'E17' = '  +2.27%  '
    'D18' = '0.109'
    'E18' = '  +1.54%  '
    'D19' = '6.78'
    'E19' = '  +2.95%  '
    'D20' = '479.17'
    'E20' = '  +0.77%  '
    'D21' = '13.77'
    'E21' = '  +4.47%  '
    'D22' = '0.688'
    'E22' = '  +2.84%  '
    'D23' = '7.59'
    'E23' = '  +8.38%  '
    'D24' = '13.50'
    'E24' = '  +12.27%  '
    'D25' = '81.15'
    'E25' = '  +0.46%  '
    'D26' = '0.999'
    'E26' = '  +0.10%  '
    'E27' = '  +3.58%  '
    'D28' = '8.18'
    'E28' = '  +6.05%  '
    'D29' = '2.06'
    'E29' = '  +7.55%  '
    'D30' = '1.00'
    'E30' = '  -0.17%  '
    'D31' = '26.13'
    'E31' = '  +2.28%  '
    'E32' = '  +2.05%  '
    'D33' = '2.48'
    'E33' = '  +6.05%  '
    'D34' = '5.61'
    'E34' = '  +0.26%  '
    'D35' = '6.16'
    'E35' = '  +5.62%  '
    'D36' = '54.97'
    'E36' = '  +0.44%  '
    'D37' = '470.71'
    'E37' = '  +2.95%  '
    'B38' = 'Hedera'
    'C38' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D38' = '0.0836'
    'E38' = '  +5.00%  '
    'B39' = 'VeChain'
    'C39' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D39' = '0.0410'
    'E39' = '  +6.47%  '
    'B40' = 'dogwifhat'
    'C40' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D40' = '2.97'
    'E40' = '  +21.93%  '
    'D41' = '2.972.63'
    'E41' = '  -5.78%  '
    'D42' = '8.26'
    'E42' = '  +2.63%  '
    'E43' = '  -1.92%  '
    'D44' = '28.24'
    'E44' = '  +7.83%  '
    'D45' = '0.259'
    'E45' = '  +7.07%  '
    'D46' = '2.15'
    'E46' = '  +9.37%  '
    'E47' = '  -0.03%  '
    'D48' = '0.112'
    'E48' = '  +3.81%  '
    'D49' = '0.0₃0527'
    'E49' = '  +6.17%  '
    'D50' = '117.13'
    'E50' = '  +2.60%  '
    'D51' = '2.07'
    'E51' = '  +3.78%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
